$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p069r_2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p069r_2</id>", 2)
